$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of dashboard data to append (rows 34-39), following the same
# layout/formatting pattern as the existing rows above them.
$data = @(
    @(842020231, "Jake Williams",  "1/2/2025", 35, "Completed"),
    @(842020232, "Tonya frison",   "1/2/2025", 35, "Completed"),
    @(842020233, "CORTNEY GRIGGS", "1/2/2025", 35, "Completed"),
    @(842020234, "Justin Roberts", "1/2/2025", 35, "Completed"),
    @(842020235, "Don Havird",     "1/2/2025", 35, "Completed"),
    @(842020236, "Scott Dieter",   "1/2/2025", 35, "Completed")
)

# Keep the same number formats used by the existing data rows so that the
# Install Date column stays text and the Total Labor Amount column stays
# currency-formatted.
$dateFormat = $ws.Cells.Item(33, 3).NumberFormat
$amountFormat = $ws.Cells.Item(33, 4).NumberFormat

$startRow = 34
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    $ws.Cells.Item($r, 3).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 3).Value = $row[2]

    $ws.Cells.Item($r, 4).NumberFormat = $amountFormat
    $ws.Cells.Item($r, 4).Value = $row[3]

    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Formula = '=HYPERLINK("https://pf.apps.projectsforce.com/project/view/8430847")'
}

# Update selection to mirror where the user finished editing.
$ws.Range("B34").Select()
